$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("model_6_0_0", 0.2024163492347106, 0.3948406239393714, -0.8065035772630404, 0.2468434048700721, 0.8826900720596313, 0.9849573969841003, 0.4652848243713379, 0.7404056191444397),
    @("model_6_0_1", 0.2240591465519546, 0.3687378973059982, -2.088985188000179, 0.06584370147039098, 0.8587379455566406, 1.027442216873169, 0.7956020832061768, 0.9183410406112671),
    @("model_6_0_3", 0.2527135858110167, 0.3869679897119189, -2.799777295980925, -0.005812749379312088, 0.827025830745697, 0.9977709650993347, 0.9786743521690369, 0.9887842535972595),
    @("model_6_0_4", 0.254914037833303, 0.3697398671084932, -2.749878594517667, -0.01476125606009449, 0.8245906233787537, 1.025811433792114, 0.9658224582672119, 0.9975812435150146),
    @("model_6_0_2", 0.2567043916014718, 0.3903084423786476, -2.568386718885416, 0.02564401552996776, 0.8226092457771301, 0.9923338294029236, 0.9190771579742432, 0.9578600525856018),
    @("model_6_0_5", 0.2683464508373961, 0.3749652719155759, -2.492830259553185, 0.02151073679989224, 0.8097248673439026, 1.017306566238403, 0.8996167778968811, 0.961923360824585),
    @("model_6_0_7", 0.2756249166049182, 0.3504862800519563, -2.376948781102824, 0.01434232118649459, 0.8016697764396667, 1.057148575782776, 0.8697702884674072, 0.9689703583717346),
    @("model_6_0_6", 0.2792865154603218, 0.3626211923249535, -2.224600140395501, 0.04376232241523248, 0.7976174354553223, 1.037397861480713, 0.8305311799049377, 0.9400485754013062),
    @("model_6_0_8", 0.2896636957891381, 0.3728271777415849, -2.204644013449982, 0.05516831745019524, 0.7861329317092896, 1.02078640460968, 0.8253912925720215, 0.9288355708122253),
    @("model_6_0_9", 0.2954430058926055, 0.3610778119283345, -2.28759023976469, 0.03464332128555125, 0.7797369360923767, 1.039909839630127, 0.8467549085617065, 0.9490131139755249),
    @("model_6_0_10", 0.2974047132925109, 0.3634082263687678, -2.270507476600951, 0.03879205430973587, 0.7775659561157227, 1.0361168384552, 0.8423551321029663, 0.944934606552124),
    @("model_6_0_11", 0.3027195546715148, 0.3450640335159366, -2.228758326712104, 0.02786033446741132, 0.7716839909553528, 1.06597375869751, 0.8316022157669067, 0.955681324005127),
    @("model_6_0_15", 0.3036970143367034, 0.3219393816740728, -2.262375052509099, 0.003446694317512211, 0.7706022262573242, 1.103611469268799, 0.8402605652809143, 0.9796815514564514),
    @("model_6_0_16", 0.3045881339673222, 0.3223176374711687, -2.244519857592166, 0.005979504747493314, 0.7696160674095154, 1.102995872497559, 0.8356617093086243, 0.9771915674209595),
    @("model_6_0_17", 0.3048910729745514, 0.3217990749685621, -2.23516462841514, 0.00667829634400452, 0.769280731678009, 1.103839874267578, 0.8332521915435791, 0.9765047430992126),
    @("model_6_0_12", 0.305194188307486, 0.3392361143428329, -2.221217036542082, 0.02368181598163199, 0.7689452767372131, 1.075459241867065, 0.8296598196029663, 0.9597890377044678),
    @("model_6_0_18", 0.3051980635596927, 0.3217114054982865, -2.228332840811606, 0.007444064099205505, 0.7689409255981445, 1.103982448577881, 0.83149254322052, 0.9757518768310547),
    @("model_6_0_13", 0.3060718651698007, 0.3336175628866314, -2.215280725977293, 0.01948907797639621, 0.7679738998413086, 1.084604024887085, 0.828130841255188, 0.9639108180999756),
    @("model_6_0_19", 0.306464747561653, 0.3205162642643946, -2.190960433447354, 0.01100410760573545, 0.7675391435623169, 1.105927705764771, 0.8218668699264526, 0.9722521305084229),
    @("model_6_0_21", 0.3066569798574272, 0.3195304975502994, -2.182372172232732, 0.01119919020850302, 0.7673264145851135, 1.107532024383545, 0.8196548819541931, 0.9720603823661804),
    @("model_6_0_14", 0.3069200327442847, 0.3328399173703311, -2.220061813018225, 0.01821774151335565, 0.767035186290741, 1.085869789123535, 0.8293622732162476, 0.9651607275009155),
    @("model_6_0_22", 0.30737564587596, 0.3212455600854491, -2.173042873273474, 0.01385226506640569, 0.7665311098098755, 1.104740619659424, 0.8172520399093628, 0.9694522619247437),
    @("model_6_0_20", 0.3074748944911292, 0.3228829100443524, -2.179764245478174, 0.0144590809238121, 0.7664212584495544, 1.102075695991516, 0.8189831972122192, 0.9688556790351868)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
    $ws.Cells.Item($row, 7).Value = $entry[6]
    $ws.Cells.Item($row, 8).Value = $entry[7]
    $ws.Cells.Item($row, 9).Value = $entry[8]
    $row++
}

# Remove rows 25 and 26 which are no longer present
$ws.Rows.Item(25).EntireRow.Delete()
$ws.Rows.Item(25).EntireRow.Delete()
